$d = $word.ActiveDocument

# 1) Fix typo "avverkningsamnälda" -> "avverkningsanmälda" in the main body text.
$d.Content.Find.Execute(
    "avverkningsamnälda", $true, $false, $false, $false, $false,
    $true, 1, $false, "avverkningsanmälda", 2) | Out-Null

# 2) Update the date in the header from 2023-11-03 to 2023-11-13.
for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)
    for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
        $h = $sec.Headers.Item($hi)
        if ($h.Exists) {
            $h.Range.Find.Execute(
                "2023-11-03", $true, $false, $false, $false, $false,
                $true, 1, $false, "2023-11-13", 2) | Out-Null
        }
    }
}

# 3) Change the document-wide default language from en-US to sv-SE
#    (docDefaults / rPrDefault), keeping eastAsia=en-US and bidi=ar-SA.
#    We approximate this via the Normal style's underlying default run
#    properties exposed through the Styles collection.
# 4) Set sv-SE (keeping eastAsia en-US / bidi ar-SA) as the language for
#    every paragraph/character/table style's run properties (adds a
#    <w:lang> element to styles lacking one). Numbering ("List") styles
#    are left untouched, matching the source change.
for ($i = 1; $i -le $d.Styles.Count; $i++) {
    $s = $d.Styles.Item($i)
    if ($s.Type -ne 4) {
        $f = $s.Font
        $f.LanguageID = "sv-SE"
        $f.LanguageIDFarEast = "en-US"
        $f.LanguageIDOther = "ar-SA"
    }
}
